$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the newly-recorded contribution rows 66-75 (previously blank
#    placeholder rows with only date-formatted C/D cells).
# ---------------------------------------------------------------------------

$rowsData = @(
    @{ Row = 66; Name = "Lunga";     Amount = 1050; Actual = 45868; MonthEnd = 45961; Type = "Contribution" },
    @{ Row = 67; Name = "Mhlengi";   Amount = 1050; Actual = 45868; MonthEnd = 45869; Type = "Contribution" },
    @{ Row = 68; Name = "Tsepo";     Amount = 1050; Actual = 45871; MonthEnd = 45869; Type = "Contribution" },
    @{ Row = 69; Name = "Piwe";      Amount = 550;  Actual = 45874; MonthEnd = 45869; Type = "Contribution" },
    @{ Row = 70; Name = "Piwe";      Amount = 300;  Actual = 45874; MonthEnd = 45869; Type = "Isipheko" },
    @{ Row = 71; Name = "Mshagmor";  Amount = 1050; Actual = 45874; MonthEnd = 45869; Type = "Contribution" },
    @{ Row = 72; Name = "Mshagmor";  Amount = 300;  Actual = 45874; MonthEnd = 45869; Type = "Isipheko" },
    @{ Row = 73; Name = "Mshagmor";  Amount = 20;   Actual = 45874; MonthEnd = 45869; Type = "Penalty" },
    @{ Row = 74; Name = "Msizi";     Amount = 1050; Actual = 45874; MonthEnd = 45869; Type = "Contribution" },
    @{ Row = 75; Name = "Msizi";     Amount = 300;  Actual = 45874; MonthEnd = 45869; Type = "Isipheko" }
)

foreach ($r in $rowsData) {
    $ws.Range("A$($r.Row)").Value = $r.Name
    $ws.Range("B$($r.Row)").Value = $r.Amount
    $ws.Range("C$($r.Row)").Value = $r.Actual
    $ws.Range("D$($r.Row)").Value = $r.MonthEnd
    $ws.Range("E$($r.Row)").Value = $r.Type
}

# ---------------------------------------------------------------------------
# 2. Rows 103/104 gain a blank, date-styled E cell (matching the style of
#    the surrounding rows, e.g. row 105). Copy formatting only, no value.
# ---------------------------------------------------------------------------

$ws.Range("C103").Copy() | Out-Null
$ws.Range("E103").PasteSpecial(-4122) | Out-Null

$ws.Range("C104").Copy() | Out-Null
$ws.Range("E104").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Drop the two trailing blank rows (371-372) that are no longer part of
#    the sheet, shrinking the used range to A1:E370.
# ---------------------------------------------------------------------------

$ws.Range("A371:E372").Delete()

# ---------------------------------------------------------------------------
# 4. Update the view: active cell / selection moves to C69 and the frozen
#    header's scroll position moves so row 54 is the first visible row.
# ---------------------------------------------------------------------------

$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C69").Select()
